# Updated cryptos list data to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.991.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.639.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("E8").Value = "  -0.74%  "
$ws.Range("E9").Value = "  -1.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.41"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0794"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.671.71"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.26%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.24"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("E14").Value = "  -0.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.24%  "
$ws.Range("E16").Value = "  -0.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.990.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("E18").Value = "  +0.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "193.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("E20").Value = "  -1.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.02%  "
$ws.Range("E22").Value = "  -2.24%  "
$ws.Range("E23").Value = "  +2.65%  "
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.78"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.55%  "
$ws.Range("B26").Value = "BinanceUSD"
$ws.Range("C26").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.82"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.46"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.62%  "
$ws.Range("E29").Value = "  -0.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0489"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.99%  "
$ws.Range("E31").Value = "  +0.35%  "
$ws.Range("E32").Value = "  -1.57%  "
$ws.Range("E33").Value = "  -1.26%  "
$ws.Range("E34").Value = "  +0.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.897"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.123.69"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.45%  "
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.532"
$ws.Range("D38").Style = "Normal"
$ws.Range("E39").Value = "  -1.18%  "
$ws.Range("E40").Value = "  -0.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.26"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.31"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.04%  "
$ws.Range("E43").Value = "  -0.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "56.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.00%  "
$ws.Range("E45").Value = "  +0.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0521"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.72"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.70%  "
$ws.Range("E49").Value = "  +0.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0940"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.17"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.55%  "
